# Diary.xlsx edit: add LeetCode problems 42 (Trapping Rain Water), 72 (Edit Distance)
# and 1143 (Longest Common Subsequence) into the sorted problem table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert row for Problem 42 "Trapping Rain Water" (goes before current row 36, A=43) ---
# Donor row 83 already has the "Bad"/"Neutral" style combination this new row needs.
$ws.Rows.Item(83).Copy()
$ws.Rows.Item(36).Insert($xlShiftDown)
$ws.Cells.Item(36, 1).Value = 42
$ws.Cells.Item(36, 2).Value = "Trapping Rain Water"
$ws.Cells.Item(36, 3).Value = ""
$ws.Cells.Item(36, 4).Value = ""

# --- Insert row for Problem 72 "Edit Distance" (goes before current row 49, A=83 after the shift above) ---
# Donor row 5 has the "Bad"/"Bad" style combination this new row needs.
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(49).Insert($xlShiftDown)
$ws.Cells.Item(49, 1).Value = 72
$ws.Cells.Item(49, 2).Value = "Edit Distance"
$ws.Cells.Item(49, 3).Value = ""
$ws.Cells.Item(49, 4).Value = ""

# --- Insert row for Problem 1143 "Longest Common Subsequence" (goes before current row 82, A=1550 after the shifts above) ---
# Donor row 19 has the "Neutral"/"Good" style combination this new row needs.
$ws.Rows.Item(19).Copy()
$ws.Rows.Item(82).Insert($xlShiftDown)
$ws.Cells.Item(82, 1).Value = 1143
$ws.Cells.Item(82, 2).Value = "Longest Common Subsequence"
$ws.Cells.Item(82, 3).Value = ""
$ws.Cells.Item(82, 4).Value = ""

$excel.CutCopyMode = $false

# The old trailing stray cell (originally C88, an empty leftover with no real
# data) shifted down to row 91 after the three inserts above - clear it out so
# the sheet's used range ends cleanly at row 90.
$ws.Rows.Item(91).Delete()

# Update the selection to match the author's final cursor position.
$ws.Range("J70").Select()
